$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the festivity dates as plain text (e.g. "2025-01-01").
# The fix simply bumps the year from 2025 to 2026 for every row, leaving
# the "nom" (B) and "ambit" (C) columns untouched.
$cellRefs  = @("A2", "A3", "A4", "A5", "A6", "A7", "A8", "A9", "A10", "A11", "A12", "A13")
$newValues = @("2026-01-01", "2026-01-06", "2026-04-03", "2026-04-06", "2026-05-01", "2026-06-24", "2026-08-15", "2026-09-11", "2026-10-12", "2026-12-08", "2026-12-25", "2026-12-26")

for ($i = 0; $i -lt $cellRefs.Length; $i++) {
    $cell = $ws.Range($cellRefs[$i])

    # Without this, assigning a "yyyy-mm-dd"-looking string makes Excel's
    # smart input parsing store it as a date serial number (with a date
    # number format) instead of the literal text the source file has.
    $cell.NumberFormat = "@"
    $cell.Value = $newValues[$i]

    # The text format was only needed to force literal-text entry; drop it
    # again so the cell keeps the workbook's original (default) styling.
    $cell.ClearFormats()
}
